# Clean up the "4U1538-52" row: the spectral type and its reference were
# edited to drop the redundant/incorrect bits:
#   ST:     "B0.2Ia" -> "B0Ia"
#   ST_ref: "Kaper (2001) Corbet et al. (2021)" -> "Kaper (2001)"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "B0Ia"
$ws.Range("C8").Value = "Kaper (2001)"

# Restore the selection to where the author last left it.
$ws.Range("C8").Select()
